$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)

$ws.Range("A6").Value = 42604.890821759262
$ws.Range("B6").Value = "Noun"
$ws.Range("C6").Value = 5876
$ws.Range("D6").Value = 3897
$ws.Range("E6").Value = 823
$ws.Range("F6").Value = 91
$ws.Range("G6").Value = 51
$ws.Range("H6").Value = 63
$ws.Range("I6").Value = 35
$ws.Range("J6").Value = 2
$ws.Range("K6").Value = 9
$ws.Range("L6").Value = 18
$ws.Range("M6").Value = 81
